$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value to a cell while forcing text storage
# (prevents Excel's autoconvert of numeric-looking strings like '1.005' into numbers),
# then restore the default 'Normal' style so no stray formatting is introduced.
function Set-TextCell($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextCell $ws.Range("D2") "28.470.24"
$ws.Range("E2").Value = "  -0.01%  "

# Row 3
Set-TextCell $ws.Range("D3") "1.825.13"

# Row 4
Set-TextCell $ws.Range("D4") "1.005"
$ws.Range("E4").Value = "  +0.40%  "

# Row 5
Set-TextCell $ws.Range("D5") "316.71"
$ws.Range("E5").Value = "  +0.47%  "

# Row 6
Set-TextCell $ws.Range("D6") "1.004"
$ws.Range("E6").Value = "  +0.28%  "

# Row 7
Set-TextCell $ws.Range("D7") "0.5163"
$ws.Range("E7").Value = "  -0.16%  "

# Row 8
Set-TextCell $ws.Range("D8") "0.3856"
$ws.Range("E8").Value = "  -1.47%  "

# Row 9
Set-TextCell $ws.Range("D9") "0.08280"
$ws.Range("E9").Value = "  +8.30%  "

# Row 10
$ws.Range("E10").Value = "  +1.24%  "

# Row 11
Set-TextCell $ws.Range("D11") "41.90"
$ws.Range("E11").Value = "  +0.04%  "

# Row 12
Set-TextCell $ws.Range("D12") "6.374"
$ws.Range("E12").Value = "  +1.12%  "

# Row 13
Set-TextCell $ws.Range("D13") "21.08"
$ws.Range("E13").Value = "  -0.24%  "

# Row 14
Set-TextCell $ws.Range("D14") "1.004"
$ws.Range("E14").Value = "  +0.31%  "

# Row 15
$ws.Range("E15").Value = "  -1.08%  "

# Row 16
Set-TextCell $ws.Range("D16") "1.823.76"
$ws.Range("E16").Value = "  -0.08%  "

# Row 17
Set-TextCell $ws.Range("D17") "94.13"
$ws.Range("E17").Value = "  +0.69%  "

# Row 18
Set-TextCell $ws.Range("D18") "0.00001121"
$ws.Range("E18").Value = "  +3.57%  "

# Row 19
Set-TextCell $ws.Range("D19") "0.06633"
$ws.Range("E19").Value = "  -0.53%  "

# Row 20
$ws.Range("E20").Value = "  +0.32%  "

# Row 21
$ws.Range("E21").Value = "  +0.28%  "

# Row 22
Set-TextCell $ws.Range("D22") "6.053"
$ws.Range("E22").Value = "  -2.33%  "

# Row 23
Set-TextCell $ws.Range("D23") "28.500.72"
$ws.Range("E23").Value = "  +0.02%  "

# Row 25
Set-TextCell $ws.Range("D25") "2.247"
$ws.Range("E25").Value = "  -0.39%  "

# Row 26
Set-TextCell $ws.Range("D26") "21.09"
$ws.Range("E26").Value = "  +2.07%  "

# Row 27
Set-TextCell $ws.Range("D27") "159.48"
$ws.Range("E27").Value = "  +1.66%  "

# Row 28
Set-TextCell $ws.Range("D28") "2.033.21"
$ws.Range("E28").Value = "  -0.11%  "

# Row 29
Set-TextCell $ws.Range("D29") "2.403"
$ws.Range("E29").Value = "  +0.01%  "

# Row 30
Set-TextCell $ws.Range("D30") "125.94"
$ws.Range("E30").Value = "  +0.60%  "

# Row 31
$ws.Range("E31").Value = "  +2.08%  "

# Row 32
Set-TextCell $ws.Range("D32") "1.093"
$ws.Range("E32").Value = "  -2.89%  "

# Row 33
Set-TextCell $ws.Range("D33") "5.729"
$ws.Range("E33").Value = "  +0.73%  "

# Row 34
Set-TextCell $ws.Range("D34") "0.07541"
$ws.Range("E34").Value = "  +7.18%  "

# Row 35
Set-TextCell $ws.Range("D35") "3.687"
$ws.Range("E35").Value = "  +0.72%  "

# Row 36
Set-TextCell $ws.Range("D36") "0.2226"
$ws.Range("E36").Value = "  -0.29%  "

# Row 37
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextCell $ws.Range("D37") "12.15"
$ws.Range("E37").Value = "  +8.21%  "

# Row 38
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell $ws.Range("D38") "0.02362"
$ws.Range("E38").Value = "  +1.52%  "

# Row 39
Set-TextCell $ws.Range("D39") "5.247"
$ws.Range("E39").Value = "  +2.04%  "

# Row 40
Set-TextCell $ws.Range("D40") "8.750"
$ws.Range("E40").Value = "  -2.55%  "

# Row 41
Set-TextCell $ws.Range("D41") "0.6388"
$ws.Range("E41").Value = "  +1.52%  "

# Row 42
Set-TextCell $ws.Range("D42") "1.187"
$ws.Range("E42").Value = "  +0.34%  "

# Row 43
Set-TextCell $ws.Range("D43") "1.396"
$ws.Range("E43").Value = "  -0.08%  "

# Row 44
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell $ws.Range("D44") "13.57"
$ws.Range("E44").Value = "  +0.83%  "

# Row 45
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextCell $ws.Range("D45") "0.6184"
$ws.Range("E45").Value = "  +4.62%  "

# Row 46
Set-TextCell $ws.Range("D46") "3.803"
$ws.Range("E46").Value = "  +2.45%  "

# Row 47
Set-TextCell $ws.Range("D47") "127.84"
$ws.Range("E47").Value = "  +2.86%  "

# Row 48
Set-TextCell $ws.Range("D48") "2.007"
$ws.Range("E48").Value = "  +1.23%  "

# Row 49
Set-TextCell $ws.Range("D49") "1.204"
$ws.Range("E49").Value = "  +0.23%  "

# Row 50
Set-TextCell $ws.Range("D50") "0.06962"
$ws.Range("E50").Value = "  +0.48%  "

# Row 51
Set-TextCell $ws.Range("D51") "1.083"
$ws.Range("E51").Value = "  +1.16%  "
